$d = $word.ActiveDocument

function Replace-ParaText($oldText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $r.Text = $newText
    }
}

# 1. Update the date in the header line
Replace-ParaText "09.07.24" "08.07.24"

# 2. Update the paper title
Replace-ParaText "Learning to (Learn at Test Time): RNNs with Expressive Hidden States" "Mixture of A Million Experts"

# 3. Update the first body paragraph
Replace-ParaText "המאמר הזה המצהיר שהוא לומד ב״זמן טסט״ משך את עיניי היום. המאמר מציע ארכיטקטורה חדשה ומעניינת לעיבוד דאטה סדרתי. בעיקרון הרשת די דומה ל-RNN מבחינת המהות אבל יש כמה הבדלים מהותיים." "המאמר של היום מציע לקחת את שיטת (Mixture of Experts(MoE לבניית ארכיטקטורות של מודלים עמוקים פופולרית במיוחד במודלי שפה. מאוד בגדול ב- MoE הרשת מורכבת מתת-רשתות (בד״כ מחלקים את שכבת ה-FFN של הטרנספורמר לכמה חלקים זרים). MoE מאומן להשתמש כל בפעם בחלק מתת-רשתות אלו (הנקראות מומחים) כאשר רשת gating רדודה יחסית באיזה מומחים צריך להשתמש כל פעם. כלומר יש לנו כן סוג של מימוש הגישה שנקראת ״lottery ticket hypothesis`" דינמי כאשר כל פעם בוחרים להריץ רק חלק מהרשת. "

# 4. Update the second body paragraph
Replace-ParaText "ֿאז מה יש לנו בארכיטקטורה הזו? בדומה ל-RNN אנו מחשבים את הייצוג עבור יחידת דאטה בזמן t (נגיד טוקן t) אבל כאן עושים זאת בשיטה שונה. לפי המאמר במקום לחשב את הייצוג עצמו אנו מחשבים את וקטור המשקלים שיאפשר לנו לחשב את ייצוגו של יחידת דאטה t. כלומר אנו מעדכנים את משקלות מודל בתנועה בהתאם לדאטה כלומר הרשת מתאפטמת ומתאימה את עצמה לדאטה שעליה היא מופעלת. זה נעשה באמצעות הזזה של המקשלים בכיוון הנגדי של הגרדיאנט של פונקציית לוס l." "כנראה שככל יש ברשת יותר מומחים בעלי אותה הארכיטקטורה וכל פעם בוחרים אותו מספר של המומחים הביצועים אמורים להשתפר אולם המחיר הוא המודל גדול יותר.המאמר מנסה לבדוק האם שווה להשתמש בהרבה מאוד במומחים רזים מאוד. המחרים מציעים לעבוד עם מיליון של מומחים של כל אחד מהם היא דל במיוחד. כמובן שכל פעם צריך לבחון את המומחים כל פעם ומכיוון שיש מיליון מומחים אז נדרש מאמץ חישובי לא קטן. המאמר מציע להשתמש בטכניקה הנקראת  product key retrieval כדי להקטין את הסיבוכיות (בגדול זה חלוקה של וקטור המפתחות (keys) לשני חלקים, ביצוע חישוב לכל אחד בנפרד ושילובם)."

# 5. Update the third body paragraph
Replace-ParaText "מה זה בעצם פונקציית l ואיך מאמנים אותה? נניח שהייצוג של איבר דאטה t מחושב על ידי פונקציית f. במקרה הזה פונקציית l יכולה להיות (למשל) נורמה של הפרש ריבוע של ייצוג דאטה z (המחושב עם f) מהדאטה עצמו. כלומר אנו מאמנים את וקטור הייצוג להיות מסוגל לשחזר (כלומר לזכור) את הדאטה עצמו x_t. כמובן שאין בזה הרבה משמעות אבל אם נאמן רשת עם קלט מורעש ונשווה את ייצוג עם הדאטה האמיתי נקבל סוג של רשת denoising שהרשת לומדת להפיק ייצוג המאפשר לזכור את הפיצ'רים המהותיים של דאטה הנחוצים לשחזור." "וגיליתי משהו מעניין במאמר הזה - יש scaling law גם ל-MoEs. אולי אסקור אותו בקרוב…"

# 6. Remove the two paragraphs describing the ttt architecture/projection approach
$pCount = $d.Paragraphs.Count
for ($i = 1; $i -le $pCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "דרך אחרת המוצעת במאמר לאמן את רשת לשחזר הטלה למימד נמוך*") {
        $pStart = $p.Range.Start
        $pNext = $d.Paragraphs.Item($i + 1)
        $pEnd = $pNext.Range.End
        $r = $d.Range($pStart, $pEnd)
        $r.Delete()
        break
    }
}

# 7. Update the arxiv link
Replace-ParaText "https://arxiv.org/pdf/2407.04620" "https://arxiv.org/abs/2407.04153"

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
